$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.992.47"
$ws.Range("E2").Value = "  +2.54%  "

$ws.Range("D3").Value = "2.342.68"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.08"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.35"
$ws.Range("E6").Value = "  +3.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.39"
$ws.Range("E10").Value = "  +4.64%  "

$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.56"
$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("E13").Value = "  +1.36%  "

$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").Value = "2.698.29"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("D17").Value = "2.330.68"
$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("D18").Value = "43.922.75"
$ws.Range("E18").Value = "  +3.14%  "

$ws.Range("E19").Value = "  +2.38%  "

$ws.Range("E20").Value = "  +1.59%  "

$ws.Range("E21").Value = "  -6.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.25"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.46"
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.61"
$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +3.42%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.55"
$ws.Range("E27").Value = "  +4.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.17"
$ws.Range("E28").Value = "  +2.93%  "

$ws.Range("E29").Value = "  -1.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.32"
$ws.Range("E30").Value = "  +5.75%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.63"
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0891"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.84"
$ws.Range("E34").Value = "  +9.97%  "

$ws.Range("E35").Value = "  +1.20%  "

$ws.Range("E36").Value = "  +4.86%  "

$ws.Range("E38").Value = "  +3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("E39").Value = "  +10.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.80"
$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("E41").Value = "  +8.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.20"
$ws.Range("E42").Value = "  +9.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  +3.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.67"
$ws.Range("E44").Value = "  +1.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("E45").Value = "  +9.30%  "

$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.29"
$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").Value = "1.672.61"
$ws.Range("E48").Value = "  -3.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.99"
$ws.Range("E49").Value = "  -2.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.97"
$ws.Range("E50").Value = "  +2.63%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.217"
$ws.Range("E51").Value = "  +15.34%  "
